# Apply "last minute updates" to the first paragraph of the document:
#   1. Give its paragraph properties a border (pBdr) matching the
#      top/left/bottom/right w:space="5" pattern already used elsewhere
#      in the document.
#   2. Increase its left indent from 120 twips (6pt) to 225 twips (11.25pt).
#   3. Replace the placeholder id text and drop the now-unneeded trailing
#      space run, merging the paragraph down to a single run.

$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)

# --- 1 & 2: paragraph formatting -------------------------------------
$pf = $p1.Range.ParagraphFormat

$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5

$pf.LeftIndent = 11.25

# --- 3: replace text, dropping the trailing space run -----------------
$r1 = $p1.Range
$fullText = $r1.Text
$trimLen = $fullText.Length - 1   # exclude the trailing paragraph mark
$idRange = $d.Range($r1.Start, $r1.Start + $trimLen)
$idRange.Text = "**ID__AFFARS_AF_PGI_5319_303__ID**"
